# Apply the commit "Update paper content and refine template" to the
# paper template document.
#
# Summary of changes:
#  1. The FirstParagraph-styled paragraph's text ("First paragraph") gets
#     the `_GoBack` bookmark re-anchored in its middle, between "First par"
#     and "agraph" (Word automatically splits the run in two around a
#     zero-length bookmark). Because a document may only have one bookmark
#     with a given name, adding the bookmark here also removes the
#     `_GoBack` bookmark that used to sit in the Title paragraph.
#  2. The (only) section becomes a "continuous" section
#     (w:type val="continuous").
#  3. The BodyText style tightens its spacing (180/180 -> 120/120 twips)
#     and becomes justified ("both").
#  4. The Bibliography style gains a hanging-indent paragraph format
#     (left=340, hanging=340 twips) and drops the space-after to 0.

$d = $word.ActiveDocument

# --- 1. Re-anchor the _GoBack bookmark inside "First paragraph" ---
# Find the paragraph using the FirstParagraph style, then split its text
# after the 9th character ("First par" | "agraph").
$firstParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "First Paragraph") {
        $firstParagraph = $p
    }
}
$splitPos = $firstParagraph.Range.Start + 9
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

# --- 2. Make the (only) section continuous ---
$d.Sections(1).PageSetup.SectionStart = 0   # wdSectionContinuous

# --- 3. BodyText style: tighter spacing + justified ---
$bodyText = $d.Styles("Body Text")
$bodyText.ParagraphFormat.SpaceBefore = 6
$bodyText.ParagraphFormat.SpaceAfter = 6
$bodyText.ParagraphFormat.Alignment = 3   # wdAlignParagraphJustify

# --- 4. Bibliography style: hanging indent + no space after ---
$bibliography = $d.Styles("Bibliography")
$bibliography.ParagraphFormat.SpaceAfter = 0
$bibliography.ParagraphFormat.LeftIndent = 17
$bibliography.ParagraphFormat.FirstLineIndent = -17
